$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Auto scs (service time, lamda_1) and auto capacity (lamda_2) were recomputed;
# every existing data row (2-51) shares the same pair of constants, now updated. ---
$newLamda1 = 33.94444444444444   # was 8.159722222222223 (time in ms)
$newLamda2 = 1.95                 # was 1.875 (auto capacity)

# --- dic_nbre_clients_poisson_2 keys/values were recomputed for the new lamda pair.
# The table also grew from 50 data rows (A1:E51) to 59 data rows (A1:E60). ---
$keys   = @(0,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,54,55,56,59,60,61,62)
$values = @(0.129,0.001,0.004,0.01,0.022,0.03,0.05,0.054,0.05,0.037,0.034,0.026,0.023,0.034,0.027,0.047,0.041,0.034,0.032,0.028,0.027,0.018,0.023,0.02,0.012,0.023,0.012,0.017,0.013,0.011,0.006,0.013,0.013,0.008,0.007,0.003,0.005,0.006,0.005,0.003,0.006,0.002,0.005,0.002,0.002,0.003,0.001,0.001,0.003,0.004,0.003,0.002,0.001,0.001,0.001,0.001,0.001,0.001,0.001)

$firstDataRow = 2
$lastDataRow  = $firstDataRow + $keys.Length - 1

for ($i = 0; $i -lt $keys.Length; $i++) {
    $r = $firstDataRow + $i

    # New rows (52-60) do not exist yet; give column A the same bold/centered/
    # bordered look as the rest of the index column before filling it in.
    if ($r -gt 51) {
        $ws.Cells.Item($r, 1).Value = $r - 2
        $ws.Cells.Item($r, 1).Font.Bold = $true
        $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
        $ws.Cells.Item($r, 1).VerticalAlignment = -4160
        $ws.Cells.Item($r, 1).Borders.LineStyle = 1
    }

    $ws.Cells.Item($r, 2).Value = $newLamda1
    $ws.Cells.Item($r, 3).Value = $newLamda2
    $ws.Cells.Item($r, 4).Value = $keys[$i]
    $ws.Cells.Item($r, 5).Value = $values[$i]
}

Write-Output "Updated rows $firstDataRow..$lastDataRow (dimension now A1:E$lastDataRow)"
